# Apply the "deleting prints, modifying plots" edit to the
# "Coupling Parameters" sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# End Year: 2030 -> 2029
$ws.Range("B4").Value = 2029

# realistic_candidate_capacities_for_future: FALSE -> TRUE
$ws.Range("B17").Value = $true

# Make this sheet active and move the selection to C17 (matches the
# saved selection in the target file).
$ws.Activate() | Out-Null
$ws.Range("C17").Select() | Out-Null
